$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 header corrections: cells that previously held placeholder
# "unnamed: ..._level_1" labels should now read "total" (matching the
# other "total" column header), and the now-unused shared strings
# "unnamed: 1_level_1" / "unnamed: 5_level_1" are dropped from the
# workbook automatically since nothing else refers to them.
$ws.Range("B2").Value = "total"
$ws.Range("F2").Value = "total"
